$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.269.97'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.516.45'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +1.26%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.30'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +5.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.49'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -1.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.513.36'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +1.49%  '

$ws.Range("E9").Value = '  -0.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.196'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +4.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.70'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +1.12%  '

$ws.Range("E12").Value = '  -2.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.43'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +0.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000279'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.078.00'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("B16").Value = 'Polkadot'

$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.37'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -4.94%  '

$ws.Range("B17").Value = 'BitcoinCash'

$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '616.74'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -9.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.512.37'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +0.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.365.37'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +1.21%  '

$ws.Range("E20").Value = '  -1.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.25'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -1.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.23'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.881'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -1.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.81'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -2.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.47'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -1.06%  '

$ws.Range("E26").Value = '  +1.93%  '

$ws.Range("E27").Value = '  +3.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("E29").Value = '  -0.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.27'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -0.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.26'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +1.42%  '

$ws.Range("E32").Value = '  -0.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.49'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -2.90%  '

$ws.Range("E34").Value = '  -0.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.94'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -4.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '570.91'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +0.97%  '

$ws.Range("B37").Value = 'Cosmos'

$ws.Range("C37").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.78'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("B38").Value = 'dogwifhat'

$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.55'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -1.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.28'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("E40").Value = '  -2.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -0.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.140'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +2.93%  '

$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.391.41'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -0.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.326'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -2.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₃0708'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +1.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.92'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -1.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.87'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -1.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.59'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +0.21%  '

$ws.Range("E50").Value = '  -2.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.57'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -0.24%  '
